$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (the existing "Jun_10" column).
# This shifts the old column C ("Jun_10" data) to column E, and pushes
# column B ("Jun_13" header) to the left unaffected.
$ws.Range("C:D").Insert()

# New header row: B1 = Jun_17 (newest), C1 = Jun_15, D1 = Jun_13 (old B1
# value, now shifted into place), E1 = Jun_10 (already shifted by Insert).
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# Fill the two new data columns (C, D) for every data row with "UN", same
# as column B, mirroring how the existing "Jun_10" column values looked
# before any rating change was recorded.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}
